$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 2150
$ws.Range("I82").Value = 500
$ws.Range("J82").Value = 3800
$ws.Range("K82").Value = 1500
$ws.Range("L82").Value = 11400
$ws.Range("M82").Value = -1094
$ws.Range("N82").Value = -12212

$ws.Range("H85").Value = 2150
$ws.Range("I85").Value = 500
$ws.Range("J85").Value = 3800
$ws.Range("K85").Value = 1500
$ws.Range("L85").Value = 11400
$ws.Range("M85").Value = -96
$ws.Range("N85").Value = -14208

$ws.Range("H106").Value = 3305.5
$ws.Range("I106").Value = 3305.5
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 3305.5
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -2674.5
$ws.Range("N106").ClearContents()

$ws.Range("H137").Value = 2643.5
$ws.Range("I137").Value = 2191.3333
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 6573.999899999999
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = -4023.999899999999
$ws.Range("N137").Value = -17100

$ws.Range("H141").Value = 5499.0586
$ws.Range("I141").Value = 2873.75
$ws.Range("J141").Value = 11799.8
$ws.Range("K141").Value = 8621.25
$ws.Range("L141").Value = 35399.39999999999
$ws.Range("M141").Value = -3441.25
$ws.Range("N141").Value = -45759.39999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 1000
$ws.Range("I13").Value = 1000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -856
$ws.Range("N13").ClearContents()

$ws.Range("H61").Value = 1982.5
$ws.Range("I61").Value = 1912.3334
$ws.Range("J61").Value = 2333.3333
$ws.Range("K61").Value = 1912.3334
$ws.Range("L61").Value = 2333.3333
$ws.Range("M61").Value = -1700.3334
$ws.Range("N61").Value = -2757.3333

$ws.Range("H88").Value = 90911190
$ws.Range("I88").Value = 1261.6666
$ws.Range("J88").Value = 200003100
$ws.Range("K88").Value = 1261.6666
$ws.Range("L88").Value = 200003100
$ws.Range("M88").Value = -855.6666
$ws.Range("N88").Value = -200003912

$ws.Range("H91").Value = 90911190
$ws.Range("I91").Value = 1261.6666
$ws.Range("J91").Value = 200003100
$ws.Range("K91").Value = 1261.6666
$ws.Range("L91").Value = 200003100
$ws.Range("M91").Value = 142.3334
$ws.Range("N91").Value = -200005908

$ws.Range("H122").Value = 1428.9736
$ws.Range("I122").Value = 1067.2
$ws.Range("J122").Value = 2785.625
$ws.Range("K122").Value = 3201.6
$ws.Range("L122").Value = 8356.875
$ws.Range("M122").Value = -751.6000000000004
$ws.Range("N122").Value = -13256.875

$ws.Range("H136").Value = 1982.5
$ws.Range("I136").Value = 1912.3334
$ws.Range("J136").Value = 2333.3333
$ws.Range("K136").Value = 5737.0002
$ws.Range("L136").Value = 6999.999899999999
$ws.Range("M136").Value = -3187.0002
$ws.Range("N136").Value = -12099.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2395.1892
$ws.Range("I20").Value = 2253.4
$ws.Range("J20").Value = 3002.8572
$ws.Range("K20").Value = 2253.4
$ws.Range("L20").Value = 3002.8572
$ws.Range("M20").Value = -2006.4
$ws.Range("N20").Value = -3496.8572

$ws.Range("H105").Value = 5861.8
$ws.Range("I105").Value = 6769.6665
$ws.Range("J105").Value = 4500
$ws.Range("K105").Value = 6769.6665
$ws.Range("L105").Value = 4500
$ws.Range("M105").Value = -5022.6665
$ws.Range("N105").Value = -7994

$ws.Range("H134").Value = 5000
$ws.Range("I134").Value = 2000
$ws.Range("J134").Value = 8000
$ws.Range("K134").Value = 6000
$ws.Range("L134").Value = 24000
$ws.Range("M134").Value = -3465
$ws.Range("N134").Value = -29070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1884.1666
$ws.Range("I5").Value = 1980
$ws.Range("J5").Value = 1405
$ws.Range("K5").Value = 5940
$ws.Range("L5").Value = 4215
$ws.Range("M5").Value = -5828
$ws.Range("N5").Value = -4439

$ws.Range("H113").Value = 740.3111
$ws.Range("I113").Value = 503.55554
$ws.Range("J113").Value = 1095.4445
$ws.Range("K113").Value = 1510.66662
$ws.Range("L113").Value = 3286.3335
$ws.Range("M113").Value = 659.33338
$ws.Range("N113").Value = -7626.333500000001

$ws.Range("H122").Value = 828.1905
$ws.Range("I122").Value = 744.75
$ws.Range("J122").Value = 847.82355
$ws.Range("K122").Value = 6702.75
$ws.Range("L122").Value = 7630.41195
$ws.Range("M122").Value = -4252.75
$ws.Range("N122").Value = -12530.41195

$ws.Range("H135").Value = 1884.1666
$ws.Range("I135").Value = 1980
$ws.Range("J135").Value = 1405
$ws.Range("K135").Value = 17820
$ws.Range("L135").Value = 12645
$ws.Range("M135").Value = -15285
$ws.Range("N135").Value = -17715

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3055.8164
$ws.Range("I80").Value = 4177.9165
$ws.Range("J80").Value = 2691.8918
$ws.Range("K80").Value = 4177.9165
$ws.Range("L80").Value = 2691.8918
$ws.Range("M80").Value = -3179.9165
$ws.Range("N80").Value = -4687.891799999999

$ws.Range("H83").Value = 3055.8164
$ws.Range("I83").Value = 4177.9165
$ws.Range("J83").Value = 2691.8918
$ws.Range("K83").Value = 20889.5825
$ws.Range("L83").Value = 13459.459
$ws.Range("M83").Value = -15897.5825
$ws.Range("N83").Value = -23443.459

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 35503.5
$ws.Range("I12").Value = 1003
$ws.Range("J12").Value = 70004
$ws.Range("K12").Value = 1003
$ws.Range("L12").Value = 70004
$ws.Range("M12").Value = -833
$ws.Range("N12").Value = -70344

$ws.Range("H82").Value = 3184.889
$ws.Range("I82").Value = 2610.2222
$ws.Range("J82").Value = 3759.5557
$ws.Range("K82").Value = 2610.2222
$ws.Range("L82").Value = 3759.5557
$ws.Range("M82").Value = -2249.2222
$ws.Range("N82").Value = -4481.5557

$ws.Range("H85").Value = 3184.889
$ws.Range("I85").Value = 2610.2222
$ws.Range("J85").Value = 3759.5557
$ws.Range("K85").Value = 2610.2222
$ws.Range("L85").Value = 3759.5557
$ws.Range("M85").Value = -1362.2222
$ws.Range("N85").Value = -6255.5557

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 200002500
$ws.Range("I81").Value = 1000000000
$ws.Range("J81").Value = 3125
$ws.Range("K81").Value = 2000000000
$ws.Range("L81").Value = 6250
$ws.Range("M81").Value = -1999998939
$ws.Range("N81").Value = -8372

$ws.Range("H84").Value = 200002500
$ws.Range("I84").Value = 1000000000
$ws.Range("J84").Value = 3125
$ws.Range("K84").Value = 10000000000
$ws.Range("L84").Value = 31250
$ws.Range("M84").Value = -9999994696
$ws.Range("N84").Value = -41858
